# The workbook tracks inventory items: column B = description, column D =
# local image file path (e.g. C:\Users\BH HUB\Desktop\Items\images\(1).jpeg).
# The edit inserts a space before the opening parenthesis in every image
# path in column D (rows 2-102), turning "...images\(1).jpeg" into
# "...images\ (1).jpeg", etc. This was done in Excel via Find & Replace
# over the whole column D (which is also why column D ends up selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column D, mirroring the user selecting the whole column before
# running Find & Replace (also reproduces the saved <selection> state).
$ws.Columns("D:D").Select()

# Insert a space before every "(" in the image-path column so
# "...images\(1).jpeg" becomes "...images\ (1).jpeg".
$range = $ws.Range("D2:D102")
$range.Replace("(", " (")
